$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 77.20339
$ws.Range("I5").Value = 77.21154
$ws.Range("J5").Value = 77.14286
$ws.Range("K5").Value = 77.21154
$ws.Range("L5").Value = 77.14286
$ws.Range("M5").Value = 37.78846
$ws.Range("N5").Value = -307.14286
$ws.Range("H6").Value = 430.4
$ws.Range("I6").Value = 430.4
$ws.Range("K6").Value = 1291.2
$ws.Range("M6").Value = -1179.2
$ws.Range("H13").Value = 2105.0952
$ws.Range("I13").Value = 300.83334
$ws.Range("J13").Value = 4510.778
$ws.Range("K13").Value = 300.83334
$ws.Range("L13").Value = 4510.778
$ws.Range("M13").Value = -131.83334
$ws.Range("N13").Value = -4848.778
$ws.Range("H43").Value = 600
$ws.Range("I43").Value = 600
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 600
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = -531
$ws.Range("N43").ClearContents()
$ws.Range("H99").Value = 3347.3076
$ws.Range("I99").Value = 213.4
$ws.Range("J99").Value = 13793.667
$ws.Range("K99").Value = 640.2
$ws.Range("L99").Value = 41381.001
$ws.Range("M99").Value = 857.8
$ws.Range("N99").Value = -44377.001
$ws.Range("H106").Value = 1252.5
$ws.Range("I106").Value = 1252.5
$ws.Range("J106").Value = 0
$ws.Range("K106").Value = 1252.5
$ws.Range("L106").Value = 0
$ws.Range("M106").Value = -621.5
$ws.Range("N106").ClearContents()
$ws.Range("H112").Value = 1106.7037
$ws.Range("I112").Value = 749.25
$ws.Range("J112").Value = 1168.8695
$ws.Range("K112").Value = 2247.75
$ws.Range("L112").Value = 3506.6085
$ws.Range("M112").Value = -1139.75
$ws.Range("N112").Value = -5722.6085
$ws.Range("H129").Value = 7063.2856
$ws.Range("I129").Value = 431.1111
$ws.Range("J129").Value = 9359.038
$ws.Range("K129").Value = 1293.3333
$ws.Range("L129").Value = 28077.114
$ws.Range("M129").Value = 3706.6667
$ws.Range("N129").Value = -38077.114
$ws.Range("H135").Value = 1369.5
$ws.Range("I135").Value = 981.6
$ws.Range("J135").Value = 3309
$ws.Range("K135").Value = 8834.4
$ws.Range("L135").Value = 29781
$ws.Range("M135").Value = -6299.4
$ws.Range("N135").Value = -34851
$ws.Range("H138").Value = 2465.7368
$ws.Range("I138").Value = 1302.9333
$ws.Range("J138").Value = 3224.087
$ws.Range("K138").Value = 3908.7999
$ws.Range("L138").Value = 9672.261
$ws.Range("M138").Value = 1231.2001
$ws.Range("N138").Value = -19952.261

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2012.5333
$ws.Range("I45").Value = 3041.8
$ws.Range("K45").Value = 3041.8
$ws.Range("M45").Value = -2664.8
$ws.Range("H94").Value = 18582.5
$ws.Range("J94").Value = 18582.5
$ws.Range("L94").Value = 18582.5
$ws.Range("N94").Value = -20384.5
$ws.Range("H97").Value = 2268.65
$ws.Range("I97").Value = 3481.9092
$ws.Range("J97").Value = 785.7778
$ws.Range("K97").Value = 3481.9092
$ws.Range("L97").Value = 785.7778
$ws.Range("M97").Value = -2985.9092
$ws.Range("N97").Value = -1777.7778
$ws.Range("H115").Value = 27336.8
$ws.Range("J115").Value = 27336.8
$ws.Range("L115").Value = 27336.8
$ws.Range("N115").Value = -30470.8

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H92").Value = 20300.334
$ws.Range("J92").Value = 20300.334
$ws.Range("L92").Value = 20300.334
$ws.Range("N92").Value = -25292.334
$ws.Range("H99").Value = 1563.2273
$ws.Range("I99").Value = 1491.5
$ws.Range("J99").Value = 1688.75
$ws.Range("K99").Value = 1491.5
$ws.Range("L99").Value = 1688.75
$ws.Range("M99").Value = 6.5
$ws.Range("N99").Value = -4684.75
$ws.Range("H114").Value = 32136.8
$ws.Range("J114").Value = 32136.8
$ws.Range("L114").Value = 32136.8
$ws.Range("N114").Value = -40814.8

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
$ws.Range("H69").Value = 4252.8
$ws.Range("I69").Value = 4252.8
$ws.Range("K69").Value = 4252.8
$ws.Range("M69").Value = -3503.8
$ws.Range("H72").Value = 4252.8
$ws.Range("I72").Value = 4252.8
$ws.Range("K72").Value = 12758.4
$ws.Range("M72").Value = -9014.400000000001
$ws.Range("H96").Value = 11975
$ws.Range("J96").Value = 11975
$ws.Range("L96").Value = 11975
$ws.Range("N96").Value = -17467
$ws.Range("H109").Value = 28799.5
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 28799.5
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 28799.5
$ws.Range("M109").ClearContents()
$ws.Range("N109").Value = -30879.5
$ws.Range("H114").Value = 27000
$ws.Range("J114").Value = 27000
$ws.Range("L114").Value = 27000
$ws.Range("N114").Value = -35678
$ws.Range("H131").Value = 56733.332
$ws.Range("J131").Value = 56733.332
$ws.Range("L131").Value = 56733.332
$ws.Range("N131").Value = -66813.33199999999
$ws.Range("H134").Value = 3725.1353
$ws.Range("I134").Value = 3702.7036
$ws.Range("J134").Value = 3785.7
$ws.Range("K134").Value = 11108.1108
$ws.Range("L134").Value = 11357.1
$ws.Range("M134").Value = -8573.110799999999
$ws.Range("N134").Value = -16427.1

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 673.1429000000001
$ws.Range("I113").Value = 397.25925
$ws.Range("J113").Value = 1011.7273
$ws.Range("K113").Value = 1191.77775
$ws.Range("L113").Value = 3035.1819
$ws.Range("M113").Value = 978.22225
$ws.Range("N113").Value = -7375.1819
$ws.Range("H122").Value = 720.3103599999999
$ws.Range("I122").Value = 320.35294
$ws.Range("J122").Value = 1286.9166
$ws.Range("K122").Value = 2883.17646
$ws.Range("L122").Value = 11582.2494
$ws.Range("M122").Value = -433.1764599999997
$ws.Range("N122").Value = -16482.2494
$ws.Range("H131").Value = 966.4375
$ws.Range("I131").Value = 370.25
$ws.Range("J131").Value = 1051.6072
$ws.Range("K131").Value = 1110.75
$ws.Range("L131").Value = 3154.8216
$ws.Range("M131").Value = 3929.25
$ws.Range("N131").Value = -13234.8216
$ws.Range("H134").Value = 1742.0834
$ws.Range("I134").Value = 2007.3334
$ws.Range("K134").Value = 6022.0002
$ws.Range("M134").Value = -952.0002000000004

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("N22").ClearContents()
$ws.Range("H28").Value = 7500
$ws.Range("J28").Value = 7500
$ws.Range("L28").Value = 7500
$ws.Range("N28").Value = -7884
$ws.Range("H33").Value = 10000
$ws.Range("J33").Value = 10000
$ws.Range("L33").Value = 10000
$ws.Range("N33").Value = -10504
$ws.Range("H95").Value = 10290.5
$ws.Range("J95").Value = 10290.5
$ws.Range("L95").Value = 10290.5
$ws.Range("N95").Value = -15782.5

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 9470.532999999999
$ws.Range("I136").Value = 30752
$ws.Range("J136").Value = 1731.8182
$ws.Range("K136").Value = 92256
$ws.Range("L136").Value = 5195.4546
$ws.Range("M136").Value = -89706
$ws.Range("N136").Value = -10295.4546
$ws.Range("H137").Value = 31450
$ws.Range("I137").Value = 13000
$ws.Range("J137").Value = 49900
$ws.Range("K137").Value = 13000
$ws.Range("L137").Value = 49900
$ws.Range("M137").Value = -7900
$ws.Range("N137").Value = -60100

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H27").Value = 18465.2
$ws.Range("J27").Value = 18465.2
$ws.Range("L27").Value = 18465.2
$ws.Range("N27").Value = -18603.2
$ws.Range("H100").Value = 29025.428
$ws.Range("I100").Value = 54522.184
$ws.Range("J100").Value = 979
$ws.Range("K100").Value = 109044.368
$ws.Range("L100").Value = 1958
$ws.Range("M100").Value = -108503.368
$ws.Range("N100").Value = -3040
$ws.Range("H107").Value = 534.64703
$ws.Range("I107").Value = 206.38461
$ws.Range("J107").Value = 1601.5
$ws.Range("K107").Value = 619.15383
$ws.Range("L107").Value = 4804.5
$ws.Range("M107").Value = 1300.84617
$ws.Range("N107").Value = -8644.5
$ws.Range("H115").Value = 19928.572
$ws.Range("J115").Value = 19928.572
$ws.Range("L115").Value = 19928.572
$ws.Range("N115").Value = -23062.572
